$wb = $excel.ActiveWorkbook

# --- Sheets ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# (shows up in the Overview summary row as well as on each language sheet's
# "Status" column)
$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# --- zh-cn sheet: handback timestamp advanced, error cleared (handback now in sync) ---
$wsZhCn.Range("K2").Value = "2016-08-31 12:23:39"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet: handback timestamp advanced, error cleared (handback now in sync) ---
$wsDeDe.Range("K2").Value = "2016-08-31 12:23:58"
$wsDeDe.Range("P2").Value = ""

# --- Column widths: re-fit now that Status / Error Detail text changed length ---
$padding = 0.8333333333333334

$wsOverview.Columns.Item(5).ColumnWidth = 29.9777050018311 - $padding
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777050018311 - $padding

$wsZhCn.Columns.Item(3).ColumnWidth  = 29.9777050018311 - $padding
$wsZhCn.Columns.Item(16).ColumnWidth = 13.7470531463623 - $padding

$wsDeDe.Columns.Item(3).ColumnWidth  = 29.9777050018311 - $padding
$wsDeDe.Columns.Item(16).ColumnWidth = 13.7470531463623 - $padding

Write-Output "Handback report regenerated"
